$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(210366942, 'Vo Nhat Tuong', 462854365, 'Kinh Doanh', 'Lam Dong'),
    @(274299112, 'Ha Son Tung', 312982000, 'Sinh Vien', 'Quang Tri'),
    @(245403375, 'Vuong Huu Toan', 461934945, 'Sinh Vien', 'Da Lat'),
    @(251262102, 'Tran Minh Tri', 732841923, 'Sinh Vien', 'TP HCM'),
    @(214587489, 'Hoang Minh Vu', 807638277, 'Sinh Vien', 'Ha Tinh '),
    @(268616154, 'Le Van', 901424984, 'Kinh Doanh', 'Nghe An'),
    @(190986546, 'Huynh Thanh Kha Tran', 650978726, 'Sinh Vien', 'Da Nang'),
    @(259937103, 'Nguyen Thi Nga', 883423258, 'Giao Vien', 'Nha Trang'),
    @(226793287, 'Thai Ba Tuyen', 596234115, 'Giao Vien', 'Quang Binh'),
    @(228568444, 'Ngo Nha Trang', 749191391, 'Kinh Doanh', 'TP HCM'),
    @(284953322, 'Nguyen Thanh Tung', 905842820, 'Sinh Vien', 'Ha Tinh '),
    @(223929938, 'Ho Minh Tuan', 924000031, 'Kinh Doanh', 'Da Lat'),
    @(216920917, 'Nguyen Van Teo', 438718377, 'Sinh Vien', 'Binh Duong'),
    @(287918095, 'Le Thanh Truong Vinh', 512555514, 'Sinh Vien', 'Nha Trang'),
    @(165024578, 'Nguyen Thi Thao', 574040645, 'Giao Vien', 'Binh Duong'),
    @(275771181, 'Le Xuan Tuan', 641379146, 'Giao Vien', 'Da Nang'),
    @(156907345, 'Huynh Luong Phuong Truc', 816762743, 'Sinh Vien', 'Dong Nai'),
    @(226851874, 'Trong Van Tu', 385950193, 'Sinh Vien', 'Lam Dong'),
    @(285206094, 'Nguyen Tho Tuan', 311602445, 'Kinh Doanh', 'Da Nang'),
    @(166952111, 'Nguyen Dinh Tuyen', 710123942, 'Sinh Vien', 'Binh Duong'),
    @(234517495, 'Bui Nhu Ngoc', 931396639, 'Giao Vien', 'TP HCM'),
    @(207232731, 'Chau Xuan Tuan', 797620758, 'Kinh Doanh', 'Nghe An'),
    @(190852217, 'Tran Thi Hang', 503885617, 'Sinh Vien', 'Quang Binh'),
    @(202321979, 'Nguyen Thi Thuy', 993768309, 'Sinh Vien', 'Quang Tri'),
    @(212029163, 'Pham Toan Trung', 353967350, 'Sinh Vien', 'Nghe An'),
    @(270869768, 'Vu Thi Hai Yen', 626876492, 'Sinh Vien', 'Nha Trang'),
    @(232958486, 'Truong Trong Tuan', 865059801, 'Sinh Vien', 'Quang Binh'),
    @(179039878, 'Ngo Quyen Trang', 967050075, 'Kinh Doanh', 'Ha Tinh '),
    @(260493543, 'Nguyen Cao Phat', 579626275, 'Kinh Doanh', 'Dong Nai'),
    @(188185163, 'La Manh Tuan', 559794230, 'Sinh Vien', 'Quang Tri')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}

$ws.Columns("F:F").ColumnWidth = 9.166666666666666

$ws.Range("B8").Select()
